$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("190824-00")

# Update B2 and C2 to be real boolean TRUE values (instead of inline string "TRUE")
$ws.Range("B2").Value = $true
$ws.Range("C2").Value = $true

# Activate the sheet and set the selection to B2:D2 with active cell B2
$ws.Activate()
$ws.Range("B2:D2").Select()
